$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we're about to touch to remain plain text,
# so values like "1.004", "0.07140", or "26.803.17" are not silently
# reinterpreted as numbers/dates (which would eat significant trailing
# zeros / thousand-style separators). D4 is intentionally excluded - its
# price value is not changing in this update. (Each contiguous block is set
# separately - a single multi-area "D2:D3,D5:D51" selector only applies the
# format to its first area.)
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = '26.803.17'
$ws.Range("D3").Value = '1.868.34'
$ws.Range("D5").Value = '292.91'
$ws.Range("D6").Value = '1.004'
$ws.Range("D7").Value = '0.5323'
$ws.Range("D8").Value = '0.3725'
$ws.Range("D9").Value = '0.07140'
$ws.Range("D10").Value = '21.35'
$ws.Range("D11").Value = '0.8849'
$ws.Range("D12").Value = '0.08131'
$ws.Range("D13").Value = '1.937.87'
$ws.Range("D14").Value = '92.16'
$ws.Range("D15").Value = '5.277'
$ws.Range("D16").Value = '1.008'
$ws.Range("D17").Value = '14.79'
$ws.Range("D18").Value = '0.000008471'
$ws.Range("D19").Value = '1.002'
$ws.Range("D20").Value = '26.753.46'
$ws.Range("D21").Value = '4.955'
$ws.Range("D22").Value = '10.62'
$ws.Range("D23").Value = '6.349'
$ws.Range("D24").Value = '2.279'
$ws.Range("D25").Value = '145.85'
$ws.Range("D26").Value = '1.738'
$ws.Range("D27").Value = '17.92'
$ws.Range("D28").Value = '113.03'
$ws.Range("D29").Value = '4.685'
$ws.Range("D30").Value = '4.607'
$ws.Range("D31").Value = '0.09084'
$ws.Range("D32").Value = '0.8024'
$ws.Range("D33").Value = '0.05015'
$ws.Range("D34").Value = '1.167'
$ws.Range("D35").Value = '2.957'
$ws.Range("D36").Value = '0.5994'
$ws.Range("D37").Value = '2.643'
$ws.Range("D38").Value = '3.189'
$ws.Range("D39").Value = '0.01933'
$ws.Range("D40").Value = '1.062'
$ws.Range("D41").Value = '0.5236'
$ws.Range("D42").Value = '6.468'
$ws.Range("D43").Value = '8.708'
$ws.Range("D44").Value = '115.52'
$ws.Range("D45").Value = '0.1484'
$ws.Range("D46").Value = '1.005'
$ws.Range("D47").Value = '10.01'
$ws.Range("D48").Value = '1.627'
$ws.Range("D49").Value = '37.15'
$ws.Range("D50").Value = '0.06047'
$ws.Range("D51").Value = '62.01'

# --- Column E (Volume(1h)) updates ---
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("E3").Value = '  -2.12%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("E5").Value = '  -4.86%  '
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  -2.51%  '
$ws.Range("E10").Value = '  -3.43%  '
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("E13").Value = '  +53.06%  '
$ws.Range("E14").Value = '  -3.86%  '
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("E18").Value = '  -2.20%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("E20").Value = '  -2.08%  '
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("E22").Value = '  -1.98%  '
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("E25").Value = '  -2.77%  '
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("E28").Value = '  -3.51%  '
$ws.Range("E29").Value = '  -2.87%  '
$ws.Range("E30").Value = '  -4.22%  '
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("E32").Value = '  -4.06%  '
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("E34").Value = '  -4.93%  '
$ws.Range("E35").Value = '  -1.28%  '
$ws.Range("E36").Value = '  +4.39%  '
$ws.Range("E37").Value = '  -1.73%  '
$ws.Range("E38").Value = '  -5.12%  '
$ws.Range("E39").Value = '  -3.81%  '
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("E41").Value = '  +6.12%  '
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("E43").Value = '  -6.45%  '
$ws.Range("E44").Value = '  -1.36%  '
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E46").Value = '  +0.39%  '
$ws.Range("E47").Value = '  -1.73%  '
$ws.Range("E48").Value = '  -0.76%  '
$ws.Range("E49").Value = '  -4.15%  '
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("E51").Value = '  -2.62%  '
